# Apply the weekly fruit/vegetable price-sheet update for the
# "Feria Lagunitas de Puerto Montt - Kiwi" sheet: rows 265-280.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 265
$ws.Cells.Item(265, 1).Value = 4
$ws.Cells.Item(265, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(265, 3).Value = "Los Lagos"
$ws.Cells.Item(265, 4).Value = 44746
$ws.Cells.Item(265, 5).Value = 10
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value = 100101
$ws.Cells.Item(265, 8).Value = "Berries"
$ws.Cells.Item(265, 9).Value = 100101007
$ws.Cells.Item(265, 10).Value = "Kiwi"
$ws.Cells.Item(265, 11).Value = "Hayward"
$ws.Cells.Item(265, 12).Value = "Primera"
$ws.Cells.Item(265, 13).Value = 200
$ws.Cells.Item(265, 14).Value = 15000
$ws.Cells.Item(265, 15).Value = 15000
$ws.Cells.Item(265, 16).Value = 15000
$ws.Cells.Item(265, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(265, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(265, 19).Value = 1000
$ws.Cells.Item(265, 20).Value = 15
$ws.Cells.Item(265, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 266
$ws.Cells.Item(266, 1).Value = 4
$ws.Cells.Item(266, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(266, 3).Value = "Los Lagos"
$ws.Cells.Item(266, 4).Value = 44746
$ws.Cells.Item(266, 5).Value = 10
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value = 100101
$ws.Cells.Item(266, 8).Value = "Berries"
$ws.Cells.Item(266, 9).Value = 100101007
$ws.Cells.Item(266, 10).Value = "Kiwi"
$ws.Cells.Item(266, 11).Value = "Hayward"
$ws.Cells.Item(266, 12).Value = "Segunda"
$ws.Cells.Item(266, 13).Value = 200
$ws.Cells.Item(266, 14).Value = 12000
$ws.Cells.Item(266, 15).Value = 12000
$ws.Cells.Item(266, 16).Value = 12000
$ws.Cells.Item(266, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(266, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(266, 19).Value = 800
$ws.Cells.Item(266, 20).Value = 15
$ws.Cells.Item(266, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 267
$ws.Cells.Item(267, 1).Value = 4
$ws.Cells.Item(267, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(267, 3).Value = "Los Lagos"
$ws.Cells.Item(267, 4).Value = 44386
$ws.Cells.Item(267, 5).Value = 10
$ws.Cells.Item(267, 6).Value = "Fruta"
$ws.Cells.Item(267, 7).Value = 100101
$ws.Cells.Item(267, 8).Value = "Berries"
$ws.Cells.Item(267, 9).Value = 100101007
$ws.Cells.Item(267, 10).Value = "Kiwi"
$ws.Cells.Item(267, 11).Value = "Hayward"
$ws.Cells.Item(267, 12).Value = "Especial"
$ws.Cells.Item(267, 13).Value = 150
$ws.Cells.Item(267, 14).Value = 15000
$ws.Cells.Item(267, 15).Value = 15000
$ws.Cells.Item(267, 16).Value = 15000
$ws.Cells.Item(267, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(267, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(267, 19).Value = 1000
$ws.Cells.Item(267, 20).Value = 15
$ws.Cells.Item(267, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 268
$ws.Cells.Item(268, 1).Value = 4
$ws.Cells.Item(268, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(268, 3).Value = "Los Lagos"
$ws.Cells.Item(268, 4).Value = 44386
$ws.Cells.Item(268, 5).Value = 10
$ws.Cells.Item(268, 6).Value = "Fruta"
$ws.Cells.Item(268, 7).Value = 100101
$ws.Cells.Item(268, 8).Value = "Berries"
$ws.Cells.Item(268, 9).Value = 100101007
$ws.Cells.Item(268, 10).Value = "Kiwi"
$ws.Cells.Item(268, 11).Value = "Hayward"
$ws.Cells.Item(268, 12).Value = "Primera"
$ws.Cells.Item(268, 13).Value = 300
$ws.Cells.Item(268, 14).Value = 13000
$ws.Cells.Item(268, 15).Value = 13500
$ws.Cells.Item(268, 16).Value = 13250
$ws.Cells.Item(268, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(268, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(268, 19).Value = 883
$ws.Cells.Item(268, 20).Value = 15
$ws.Cells.Item(268, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 269
$ws.Cells.Item(269, 1).Value = 4
$ws.Cells.Item(269, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(269, 3).Value = "Los Lagos"
$ws.Cells.Item(269, 4).Value = 44690
$ws.Cells.Item(269, 5).Value = 10
$ws.Cells.Item(269, 6).Value = "Fruta"
$ws.Cells.Item(269, 7).Value = 100101
$ws.Cells.Item(269, 8).Value = "Berries"
$ws.Cells.Item(269, 9).Value = 100101007
$ws.Cells.Item(269, 10).Value = "Kiwi"
$ws.Cells.Item(269, 11).Value = "Hayward"
$ws.Cells.Item(269, 12).Value = "Especial"
$ws.Cells.Item(269, 13).Value = 200
$ws.Cells.Item(269, 14).Value = 20000
$ws.Cells.Item(269, 15).Value = 20000
$ws.Cells.Item(269, 16).Value = 20000
$ws.Cells.Item(269, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(269, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(269, 19).Value = 1333
$ws.Cells.Item(269, 20).Value = 15
$ws.Cells.Item(269, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 270
$ws.Cells.Item(270, 1).Value = 4
$ws.Cells.Item(270, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(270, 3).Value = "Los Lagos"
$ws.Cells.Item(270, 4).Value = 44690
$ws.Cells.Item(270, 5).Value = 10
$ws.Cells.Item(270, 6).Value = "Fruta"
$ws.Cells.Item(270, 7).Value = 100101
$ws.Cells.Item(270, 8).Value = "Berries"
$ws.Cells.Item(270, 9).Value = 100101007
$ws.Cells.Item(270, 10).Value = "Kiwi"
$ws.Cells.Item(270, 11).Value = "Hayward"
$ws.Cells.Item(270, 12).Value = "Primera"
$ws.Cells.Item(270, 13).Value = 200
$ws.Cells.Item(270, 14).Value = 18000
$ws.Cells.Item(270, 15).Value = 18000
$ws.Cells.Item(270, 16).Value = 18000
$ws.Cells.Item(270, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(270, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(270, 19).Value = 1200
$ws.Cells.Item(270, 20).Value = 15
$ws.Cells.Item(270, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 271
$ws.Cells.Item(271, 1).Value = 4
$ws.Cells.Item(271, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(271, 3).Value = "Los Lagos"
$ws.Cells.Item(271, 4).Value = 44690
$ws.Cells.Item(271, 5).Value = 10
$ws.Cells.Item(271, 6).Value = "Fruta"
$ws.Cells.Item(271, 7).Value = 100101
$ws.Cells.Item(271, 8).Value = "Berries"
$ws.Cells.Item(271, 9).Value = 100101007
$ws.Cells.Item(271, 10).Value = "Kiwi"
$ws.Cells.Item(271, 11).Value = "Hayward"
$ws.Cells.Item(271, 12).Value = "Segunda"
$ws.Cells.Item(271, 13).Value = 200
$ws.Cells.Item(271, 14).Value = 16000
$ws.Cells.Item(271, 15).Value = 16000
$ws.Cells.Item(271, 16).Value = 16000
$ws.Cells.Item(271, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(271, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(271, 19).Value = 1067
$ws.Cells.Item(271, 20).Value = 15
$ws.Cells.Item(271, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 272
$ws.Cells.Item(272, 1).Value = 4
$ws.Cells.Item(272, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(272, 3).Value = "Los Lagos"
$ws.Cells.Item(272, 4).Value = 44672
$ws.Cells.Item(272, 5).Value = 10
$ws.Cells.Item(272, 6).Value = "Fruta"
$ws.Cells.Item(272, 7).Value = 100101
$ws.Cells.Item(272, 8).Value = "Berries"
$ws.Cells.Item(272, 9).Value = 100101007
$ws.Cells.Item(272, 10).Value = "Kiwi"
$ws.Cells.Item(272, 11).Value = "Hayward"
$ws.Cells.Item(272, 12).Value = "Primera"
$ws.Cells.Item(272, 13).Value = 400
$ws.Cells.Item(272, 14).Value = 17000
$ws.Cells.Item(272, 15).Value = 18000
$ws.Cells.Item(272, 16).Value = 17500
$ws.Cells.Item(272, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(272, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(272, 19).Value = 1167
$ws.Cells.Item(272, 20).Value = 15
$ws.Cells.Item(272, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 273
$ws.Cells.Item(273, 1).Value = 4
$ws.Cells.Item(273, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(273, 3).Value = "Los Lagos"
$ws.Cells.Item(273, 4).Value = 44672
$ws.Cells.Item(273, 5).Value = 10
$ws.Cells.Item(273, 6).Value = "Fruta"
$ws.Cells.Item(273, 7).Value = 100101
$ws.Cells.Item(273, 8).Value = "Berries"
$ws.Cells.Item(273, 9).Value = 100101007
$ws.Cells.Item(273, 10).Value = "Kiwi"
$ws.Cells.Item(273, 11).Value = "Hayward"
$ws.Cells.Item(273, 12).Value = "Segunda"
$ws.Cells.Item(273, 13).Value = 200
$ws.Cells.Item(273, 14).Value = 15000
$ws.Cells.Item(273, 15).Value = 15000
$ws.Cells.Item(273, 16).Value = 15000
$ws.Cells.Item(273, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(273, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(273, 19).Value = 1000
$ws.Cells.Item(273, 20).Value = 15
$ws.Cells.Item(273, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 274
$ws.Cells.Item(274, 1).Value = 4
$ws.Cells.Item(274, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(274, 3).Value = "Los Lagos"
$ws.Cells.Item(274, 4).Value = 44344
$ws.Cells.Item(274, 5).Value = 10
$ws.Cells.Item(274, 6).Value = "Fruta"
$ws.Cells.Item(274, 7).Value = 100101
$ws.Cells.Item(274, 8).Value = "Berries"
$ws.Cells.Item(274, 9).Value = 100101007
$ws.Cells.Item(274, 10).Value = "Kiwi"
$ws.Cells.Item(274, 11).Value = "Hayward"
$ws.Cells.Item(274, 12).Value = "Especial"
$ws.Cells.Item(274, 13).Value = 200
$ws.Cells.Item(274, 14).Value = 17000
$ws.Cells.Item(274, 15).Value = 17000
$ws.Cells.Item(274, 16).Value = 17000
$ws.Cells.Item(274, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(274, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(274, 19).Value = 1133
$ws.Cells.Item(274, 20).Value = 15
$ws.Cells.Item(274, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 275
$ws.Cells.Item(275, 1).Value = 4
$ws.Cells.Item(275, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(275, 3).Value = "Los Lagos"
$ws.Cells.Item(275, 4).Value = 44421
$ws.Cells.Item(275, 5).Value = 10
$ws.Cells.Item(275, 6).Value = "Fruta"
$ws.Cells.Item(275, 7).Value = 100101
$ws.Cells.Item(275, 8).Value = "Berries"
$ws.Cells.Item(275, 9).Value = 100101007
$ws.Cells.Item(275, 10).Value = "Kiwi"
$ws.Cells.Item(275, 11).Value = "Hayward"
$ws.Cells.Item(275, 12).Value = "Especial"
$ws.Cells.Item(275, 13).Value = 400
$ws.Cells.Item(275, 14).Value = 21000
$ws.Cells.Item(275, 15).Value = 21000
$ws.Cells.Item(275, 16).Value = 21000
$ws.Cells.Item(275, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(275, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(275, 19).Value = 1400
$ws.Cells.Item(275, 20).Value = 15
$ws.Cells.Item(275, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 276
$ws.Cells.Item(276, 1).Value = 4
$ws.Cells.Item(276, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(276, 3).Value = "Los Lagos"
$ws.Cells.Item(276, 4).Value = 44421
$ws.Cells.Item(276, 5).Value = 10
$ws.Cells.Item(276, 6).Value = "Fruta"
$ws.Cells.Item(276, 7).Value = 100101
$ws.Cells.Item(276, 8).Value = "Berries"
$ws.Cells.Item(276, 9).Value = 100101007
$ws.Cells.Item(276, 10).Value = "Kiwi"
$ws.Cells.Item(276, 11).Value = "Hayward"
$ws.Cells.Item(276, 12).Value = "Primera"
$ws.Cells.Item(276, 13).Value = 200
$ws.Cells.Item(276, 14).Value = 14000
$ws.Cells.Item(276, 15).Value = 14000
$ws.Cells.Item(276, 16).Value = 14000
$ws.Cells.Item(276, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(276, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(276, 19).Value = 933
$ws.Cells.Item(276, 20).Value = 15
$ws.Cells.Item(276, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 277
$ws.Cells.Item(277, 1).Value = 4
$ws.Cells.Item(277, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(277, 3).Value = "Los Lagos"
$ws.Cells.Item(277, 4).Value = 44442
$ws.Cells.Item(277, 5).Value = 10
$ws.Cells.Item(277, 6).Value = "Fruta"
$ws.Cells.Item(277, 7).Value = 100101
$ws.Cells.Item(277, 8).Value = "Berries"
$ws.Cells.Item(277, 9).Value = 100101007
$ws.Cells.Item(277, 10).Value = "Kiwi"
$ws.Cells.Item(277, 11).Value = "Hayward"
$ws.Cells.Item(277, 12).Value = "Especial"
$ws.Cells.Item(277, 13).Value = 300
$ws.Cells.Item(277, 14).Value = 21000
$ws.Cells.Item(277, 15).Value = 21000
$ws.Cells.Item(277, 16).Value = 21000
$ws.Cells.Item(277, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(277, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(277, 19).Value = 1400
$ws.Cells.Item(277, 20).Value = 15
$ws.Cells.Item(277, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 278
$ws.Cells.Item(278, 1).Value = 4
$ws.Cells.Item(278, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value = "Los Lagos"
$ws.Cells.Item(278, 4).Value = 44442
$ws.Cells.Item(278, 5).Value = 10
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100101
$ws.Cells.Item(278, 8).Value = "Berries"
$ws.Cells.Item(278, 9).Value = 100101007
$ws.Cells.Item(278, 10).Value = "Kiwi"
$ws.Cells.Item(278, 11).Value = "Hayward"
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 150
$ws.Cells.Item(278, 14).Value = 14000
$ws.Cells.Item(278, 15).Value = 14000
$ws.Cells.Item(278, 16).Value = 14000
$ws.Cells.Item(278, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(278, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(278, 19).Value = 933
$ws.Cells.Item(278, 20).Value = 15
$ws.Cells.Item(278, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 279
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value = 44519
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100101
$ws.Cells.Item(279, 8).Value = "Berries"
$ws.Cells.Item(279, 9).Value = 100101007
$ws.Cells.Item(279, 10).Value = "Kiwi"
$ws.Cells.Item(279, 11).Value = "Hayward"
$ws.Cells.Item(279, 12).Value = "Especial"
$ws.Cells.Item(279, 13).Value = 200
$ws.Cells.Item(279, 14).Value = 22000
$ws.Cells.Item(279, 15).Value = 22000
$ws.Cells.Item(279, 16).Value = 22000
$ws.Cells.Item(279, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(279, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(279, 19).Value = 1467
$ws.Cells.Item(279, 20).Value = 15
$ws.Cells.Item(279, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 280
$ws.Cells.Item(280, 1).Value = 4
$ws.Cells.Item(280, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(280, 3).Value = "Los Lagos"
$ws.Cells.Item(280, 4).Value = 44519
$ws.Cells.Item(280, 5).Value = 10
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100101
$ws.Cells.Item(280, 8).Value = "Berries"
$ws.Cells.Item(280, 9).Value = 100101007
$ws.Cells.Item(280, 10).Value = "Kiwi"
$ws.Cells.Item(280, 11).Value = "Hayward"
$ws.Cells.Item(280, 12).Value = "Primera"
$ws.Cells.Item(280, 13).Value = 500
$ws.Cells.Item(280, 14).Value = 15000
$ws.Cells.Item(280, 15).Value = 16000
$ws.Cells.Item(280, 16).Value = 15500
$ws.Cells.Item(280, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(280, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(280, 19).Value = 1033
$ws.Cells.Item(280, 20).Value = 15
$ws.Cells.Item(280, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

"Updated rows 265-280 on " + $ws.Name
